$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.046200752258301
$ws.Range("B1").Value = 2.080154895782471
$ws.Range("C1").Value = 2.244231939315796
$ws.Range("D1").Value = 3.057664155960083
$ws.Range("E1").Value = 2.913592338562012
